$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grupos")

# Update the active cell selection shown in the sheet view
$ws.Range("J6").Select()

# Row 2
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 3

# Row 3
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 3

# Row 4
$ws.Range("E4").Value = 20
$ws.Range("F4").Value = -17
$ws.Range("G4").Value = 3

# Row 5
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 3
